$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 93
$ws.Range("I2").Value = 302
$ws.Range("J2").Value = 1153
$ws.Range("L2").Value = 322
$ws.Range("M2").Value = 29
$ws.Range("N2").Value = 202
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 16
$ws.Range("S2").Value = 125
$ws.Range("T2").Value = 218
$ws.Range("V2").Value = 1738
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1793
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 14
